$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws 'D2' '28.465.76'
$ws.Range('E2').Value = '  -0.94%  '

Set-TextValue $ws 'D3' '1.862.33'
$ws.Range('E3').Value = '  -0.56%  '

$ws.Range('E4').Value = '  +0.21%  '

Set-TextValue $ws 'D5' '324.82'
$ws.Range('E5').Value = '  -0.70%  '

Set-TextValue $ws 'D6' '1.006'
$ws.Range('E6').Value = '  +0.05%  '

Set-TextValue $ws 'D7' '0.4553'
$ws.Range('E7').Value = '  -2.34%  '

Set-TextValue $ws 'D8' '0.3827'
$ws.Range('E8').Value = '  -2.36%  '

Set-TextValue $ws 'D9' '0.07811'
$ws.Range('E9').Value = '  -1.34%  '

Set-TextValue $ws 'D10' '0.9847'
$ws.Range('E10').Value = '  +0.84%  '

Set-TextValue $ws 'D11' '21.47'
$ws.Range('E11').Value = '  -4.02%  '

Set-TextValue $ws 'D12' '1.860.92'
$ws.Range('E12').Value = '  -0.60%  '

Set-TextValue $ws 'D13' '6.896'
$ws.Range('E13').Value = '  -1.00%  '

Set-TextValue $ws 'D14' '5.630'
$ws.Range('E14').Value = '  -1.98%  '

Set-TextValue $ws 'D15' '0.06913'
$ws.Range('E15').Value = '  -0.40%  '

Set-TextValue $ws 'D16' '1.007'
$ws.Range('E16').Value = '  +0.16%  '

Set-TextValue $ws 'D17' '86.46'
$ws.Range('E17').Value = '  -2.82%  '

Set-TextValue $ws 'D18' '0.000009919'
$ws.Range('E18').Value = '  -1.30%  '

Set-TextValue $ws 'D19' '16.69'
$ws.Range('E19').Value = '  -1.77%  '

$ws.Range('E20').Value = '  +0.01%  '

Set-TextValue $ws 'D21' '28.467.77'
$ws.Range('E21').Value = '  -0.78%  '

$ws.Range('E22').Value = '  -1.81%  '

Set-TextValue $ws 'D23' '10.87'
$ws.Range('E23').Value = '  -2.25%  '

Set-TextValue $ws 'D24' '2.086'
$ws.Range('E24').Value = '  -1.89%  '

Set-TextValue $ws 'D25' '2.063.60'
$ws.Range('E25').Value = '  -1.40%  '

Set-TextValue $ws 'D26' '153.20'
$ws.Range('E26').Value = '  -1.48%  '

Set-TextValue $ws 'D27' '19.08'
$ws.Range('E27').Value = '  -1.38%  '

Set-TextValue $ws 'D28' '5.645'
$ws.Range('E28').Value = '  -2.53%  '

Set-TextValue $ws 'D29' '117.29'
$ws.Range('E29').Value = '  -1.94%  '

Set-TextValue $ws 'D30' '1.890'
$ws.Range('E30').Value = '  -5.52%  '

Set-TextValue $ws 'D31' '0.09264'
$ws.Range('E31').Value = '  -1.08%  '

Set-TextValue $ws 'D32' '0.9033'
$ws.Range('E32').Value = '  -4.17%  '

Set-TextValue $ws 'D33' '5.269'
$ws.Range('E33').Value = '  -1.27%  '

$ws.Range('E34').Value = '  -2.53%  '

Set-TextValue $ws 'D35' '3.285'
$ws.Range('E35').Value = '  -1.93%  '

Set-TextValue $ws 'D36' '0.05666'
$ws.Range('E36').Value = '  -3.15%  '

Set-TextValue $ws 'D37' '1.147'
$ws.Range('E37').Value = '  -0.94%  '

Set-TextValue $ws 'D38' '0.02038'
$ws.Range('E38').Value = '  -3.95%  '

Set-TextValue $ws 'D39' '7.621'
$ws.Range('E39').Value = '  -3.16%  '

Set-TextValue $ws 'D40' '0.5546'
$ws.Range('E40').Value = '  -2.21%  '

Set-TextValue $ws 'D41' '0.1761'
$ws.Range('E41').Value = '  -1.15%  '

Set-TextValue $ws 'D42' '9.587'
$ws.Range('E42').Value = '  -4.11%  '

Set-TextValue $ws 'D43' '0.07126'
$ws.Range('E43').Value = '  -3.20%  '

Set-TextValue $ws 'D44' '11.51'
$ws.Range('E44').Value = '  -1.33%  '

Set-TextValue $ws 'D45' '0.5227'
$ws.Range('E45').Value = '  -2.02%  '

Set-TextValue $ws 'D46' '1.120'
$ws.Range('E46').Value = '  -1.72%  '

Set-TextValue $ws 'D47' '2.095'
$ws.Range('E47').Value = '  -5.75%  '

$ws.Range('E48').Value = '  -2.88%  '

Set-TextValue $ws 'D49' '111.58'
$ws.Range('E49').Value = '  -2.27%  '

Set-TextValue $ws 'D50' '2.430'
$ws.Range('E50').Value = '  +3.13%  '

Set-TextValue $ws 'D51' '1.006'
$ws.Range('E51').Value = '  -0.04%  '
